$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.802.79"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "2.400.02"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("D5").Value = "561.54"
$ws.Range("E5").Value = "  -1.60%  "
$ws.Range("D6").Value = "141.98"
$ws.Range("E6").Value = "  +1.46%  "
$ws.Range("D8").Value = "0.538"
$ws.Range("E8").Value = "  +2.14%  "
$ws.Range("D9").Value = "2.406.85"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("E12").Value = "  +1.56%  "
$ws.Range("D14").Value = "26.13"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.843.44"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0000169"
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("D17").Value = "60.682.01"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").Value = "2.411.49"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").Value = "8.17"
$ws.Range("E19").Value = "  +7.65%  "
$ws.Range("D20").Value = "10.67"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").Value = "323.82"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").Value = "6.04"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").Value = "1.84"
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("D26").Value = "64.63"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("D27").Value = "572.65"
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("D28").Value = "8.03"
$ws.Range("E29").Value = "  -1.23%  "
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("D31").Value = "8.07"
$ws.Range("E31").Value = "  +2.06%  "
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("E33").Value = "  -2.16%  "
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("E36").Value = "  +3.18%  "
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("D40").Value = "18.28"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("D42").Value = "2.55"
$ws.Range("E42").Value = "  +7.98%  "
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "1.68"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "41.67"
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("D46").Value = "0.0₆0278"
$ws.Range("E46").Value = "  -4.22%  "
$ws.Range("D47").Value = "142.32"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").Value = "3.52"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("E50").Value = "  +1.04%  "
$ws.Range("D51").Value = "19.34"
$ws.Range("E51").Value = "  -1.30%  "
